$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.056.92"
$ws.Range("E2").Value = "  -1.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.382.35"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.02%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.32%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.381.69"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -1.31%  "

# Row 10 - Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.87%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -3.10%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.63%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.962.85"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.64%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.38%  "

# Row 16 - ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.33%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.385.48"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "61.227.01"

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -2.76%  "

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.05%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.63%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.62%  "

# Row 23 - WrappedeETH
$ws.Range("D23").Value = "3.519.96"
$ws.Range("E23").Value = "  -0.28%  "

# Row 24 - Polygon
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.549"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.54%  "

# Row 26 - now Litecoin (was PEPE)
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.42%  "

# Row 27 - now PEPE (was Litecoin)
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.26%  "

# Row 28 - now Fetch.AI (was Kaspa)
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.35%  "

# Row 29 - now Kaspa (was Fetch.AI)
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.179"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.88%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.14%  "

# Row 31 - RenderToken
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.93%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.78%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -1.67%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.03%  "

# Row 35 - EthereumClassic
$ws.Range("E35").Value = "  -0.19%  "

# Row 36 - NEARProtocol
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.14"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.21%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -1.27%  "

# Row 38 - Aptos
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.85%  "

# Row 39 - Monero
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.63%  "

# Row 40 - Hedera
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0761"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.33%  "

# Row 41 - now FirstDigitalUSD (was EnergySwap)
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.08%  "

# Row 42 - now EnergySwap (was FirstDigitalUSD)
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.37"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.73%  "

# Row 43 - Mantle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.777"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.77%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -1.93%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  -2.14%  "

# Row 46 - ONDO
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.64%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.553.28"
$ws.Range("E47").Value = "  +8.81%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  -1.29%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.13%  "

# Row 50 - LidoDAOToken
$ws.Range("E50").Value = "  +3.65%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -1.30%  "
